# Translate Swedish "Rad X, Kol Y" labels to English "Row X, Col Y"
# in the worksheet, leaving all other cell contents (and formatting)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Row 1, Col 1"
$ws.Range("C1").Value = "Row 1, Col 3"
$ws.Range("D1").Value = "Row 1, Col 4"
$ws.Range("E1").Value = "Row 1, Col 5"
$ws.Range("F1").Value = "Row 1, Col 6"

$ws.Range("A2").Value = "Row 2, Col 1"
$ws.Range("C2").Value = "Row 2, Col 3"
$ws.Range("D2").Value = "Row 2, Col 4"
$ws.Range("F2").Value = "Row 2, Col 6"

$ws.Range("A3").Value = "Row 3, Col 1"
$ws.Range("C3").Value = "Row 3, Col 3"
$ws.Range("D3").Value = "Row 3, Col 4"
$ws.Range("E3").Value = "Row 3, Col 5"

$ws.Range("A4").Value = "Row 4, Col 1"
$ws.Range("D4").Value = "Row 4, Col 4"
$ws.Range("E4").Value = "Row 4, Col 5"
$ws.Range("F4").Value = "Row 4, Col 6"

$ws.Range("A5").Value = "Row 5, Col 1"
$ws.Range("C5").Value = "Row 5, Col 3"
$ws.Range("D5").Value = "Row 5, Col 4"

$ws.Range("C6").Value = "Row 6, Col 3"
$ws.Range("D6").Value = "Row 6, Col 4"

$ws.Range("A7").Value = "Row 7, Col 1"
$ws.Range("D7").Value = "Row 7, Col 4"
$ws.Range("E7").Value = "Row 7, Col 5"
$ws.Range("F7").Value = "Row 7, Col 6"
